# Update workbook metadata: reclassify sector-descripcion, sexo and
# mes-y-ano columns from "measure" to "dimension", and attach their
# corresponding external mapping files in a new row 6.
# (Actualización de datos obtenidos el 6 de abril de 2016)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: iaest-measure:* -> iaest-dimension:* for sector-descripcion, sexo, mes-y-ano
$ws.Range("A3").Value = "iaest-dimension:sector-descripcion"
$ws.Range("F3").Value = "iaest-dimension:sector-descripcion"
$ws.Range("K3").Value = "iaest-dimension:sexo"
$ws.Range("L3").Value = "iaest-dimension:mes-y-ano"

# Row 4: medida -> dim for the same columns
$ws.Range("A4").Value = "dim"
$ws.Range("F4").Value = "dim"
$ws.Range("K4").Value = "dim"
$ws.Range("L4").Value = "dim"

# Row 5: xsd:string -> skos:Concept for sector-descripcion / sexo columns
# (mes-y-ano, column L, keeps its xsd:string datatype)
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("F5").Value = "skos:Concept"
$ws.Range("K5").Value = "skos:Concept"

# New row 6: external mapping file references.
# Copy the formatting from row 5 first so the new cells share the same
# style as the rest of the table, then fill in the values.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("F5").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("K5").Copy()
$ws.Range("K6").PasteSpecial(-4122)

$ws.Range("A6").Value = "mapping-sector-descripcion.xlsx"
$ws.Range("F6").Value = "mapping-sector-descripcion.xlsx"
$ws.Range("K6").Value = "mapping-sexo.xlsx"
